$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 04:05"

# Update per-country statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)
# Row 5
$ws.Range("B5").Value = 499966
$ws.Range("C5").Value = 1526
$ws.Range("D5").Value = 205371
$ws.Range("E5").Value = 265746
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 15
$ws.Range("H5").Value = 28849

# Row 51
$ws.Range("B51").Value = 11468
$ws.Range("C51").Value = 27
$ws.Range("D51").Value = 10405
$ws.Range("E51").Value = 793
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 270

# Row 57
$ws.Range("B57").Value = 9592
$ws.Range("C57").Value = 861
$ws.Range("D57").Value = 889
$ws.Range("E57").Value = 8393
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 10
$ws.Range("H57").Value = 310

# Row 58
$ws.Range("B58").Value = 9267
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 5549
$ws.Range("E58").Value = 3072
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 646

# Row 59
$ws.Range("B59").Value = 9230
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 6546
$ws.Range("E59").Value = 2365
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 319

# Row 60
$ws.Range("B60").Value = 8927
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 3317
$ws.Range("E60").Value = 5483
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 127

# Row 66
$ws.Range("B66").Value = 7192
$ws.Range("C66").Value = 7
$ws.Range("D66").Value = 6614
$ws.Range("E66").Value = 475
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 103

# Row 71
$ws.Range("B71").Value = 5094
$ws.Range("C71").Value = 208
$ws.Range("D71").Value = 536
$ws.Range("E71").Value = 4357
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 201

# Row 73
$ws.Range("B73").Value = 4739
$ws.Range("C73").Value = 132
$ws.Range("D73").Value = 706
$ws.Range("E73").Value = 3931
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 12
$ws.Range("H73").Value = 102

# Row 88
$ws.Range("D88").Value = 1031
$ws.Range("E88").Value = 1318
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 46

# Row 94
$ws.Range("C94").Value = 281
$ws.Range("D94").Value = 24
$ws.Range("E94").Value = 1800
$ws.Range("G94").Value = 6
$ws.Range("H94").Value = 41

# Row 95
$ws.Range("B95").Value = 1865
$ws.Range("D95").Value = 1622
$ws.Range("E95").Value = 176
$ws.Range("H95").Value = 67

# Row 96
$ws.Range("B96").Value = 1806
$ws.Range("D96").Value = 1794
$ws.Range("E96").Value = 2
$ws.Range("H96").Value = 10

# Row 97
$ws.Range("B97").Value = 1722
$ws.Range("D97").Value = 1113
$ws.Range("E97").Value = 593
$ws.Range("H97").Value = 16

# Row 98
$ws.Range("B98").Value = 1699
$ws.Range("D98").Value = 1385
$ws.Range("E98").Value = 293
$ws.Range("H98").Value = 21

# Row 99
$ws.Range("B99").Value = 1672
$ws.Range("D99").Value = 406
$ws.Range("E99").Value = 1261
$ws.Range("H99").Value = 5

# Row 100
$ws.Range("B100").Value = 1670
$ws.Range("D100").Value = 1229
$ws.Range("E100").Value = 371
$ws.Range("H100").Value = 70

# Row 101
$ws.Range("B101").Value = 1620
$ws.Range("C101").Value = 7
$ws.Range("D101").Value = 781
$ws.Range("E101").Value = 829
$ws.Range("H101").Value = 10

# Row 172
$ws.Range("B172").Value = 117
$ws.Range("C172").Value = 1
$ws.Range("E172").Value = 1

# Row 213
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

# Row 214
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
